$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.106.66"
$ws.Range("E2").Value = "  +6.28%  "
$ws.Range("D3").Value = "3.114.84"
$ws.Range("E3").Value = "  +3.90%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.69"
$ws.Range("E5").Value = "  +4.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.03"
$ws.Range("E6").Value = "  +3.69%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "3.103.64"
$ws.Range("E8").Value = "  +3.94%  "
$ws.Range("E9").Value = "  +2.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.147"
$ws.Range("E10").Value = "  +10.44%  "
$ws.Range("E11").Value = "  +9.56%  "
$ws.Range("E12").Value = "  +2.03%  "
$ws.Range("E13").Value = "  +6.30%  "
$ws.Range("E14").Value = "  +5.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.124"
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("D16").Value = "3.628.15"
$ws.Range("E16").Value = "  +3.84%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.26"
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("D18").Value = "63.031.26"
$ws.Range("E18").Value = "  +6.19%  "
$ws.Range("D19").Value = "3.111.20"
$ws.Range("E19").Value = "  +3.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "454.90"
$ws.Range("E20").Value = "  +5.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.11"
$ws.Range("E21").Value = "  +3.54%  "
$ws.Range("E22").Value = "  +1.69%  "
$ws.Range("E23").Value = "  +6.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.65"
$ws.Range("E24").Value = "  +0.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.14"
$ws.Range("E25").Value = "  +2.23%  "
$ws.Range("E27").Value = "  +1.23%  "
$ws.Range("E28").Value = "  +6.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.32"
$ws.Range("E29").Value = "  +5.31%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.89"
$ws.Range("E31").Value = "  +12.31%  "
$ws.Range("E32").Value = "  +10.84%  "
$ws.Range("E33").Value = "  +5.19%  "
$ws.Range("D34").Value = "0.0₃0824"
$ws.Range("E34").Value = "  +7.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.37"
$ws.Range("E35").Value = "  +11.80%  "
$ws.Range("E36").Value = "  +3.70%  "
$ws.Range("E37").Value = "  +1.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.11"
$ws.Range("E38").Value = "  +12.80%  "
$ws.Range("E39").Value = "  +4.15%  "
$ws.Range("E40").Value = "  +1.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "428.48"
$ws.Range("E41").Value = "  +5.40%  "
$ws.Range("D42").Value = "2.971.39"
$ws.Range("E42").Value = "  +6.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0374"
$ws.Range("E43").Value = "  +5.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.112"
$ws.Range("E44").Value = "  +3.86%  "
$ws.Range("E45").Value = "  +8.94%  "
$ws.Range("E46").Value = "  +7.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "125.08"
$ws.Range("E47").Value = "  +1.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.75"
$ws.Range("E49").Value = "  -0.52%  "
$ws.Range("E50").Value = "  +1.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.92"
$ws.Range("E51").Value = "  +5.71%  "
